# AllSoundList.xlsx - "fix - stage bgm 변경"
#
# 1) All sound file paths in column B switch from forward-slash relative
#    paths ("../../Resource/sounds/...") to backslash-style paths
#    ("\..\Resource\sounds\...").
# 2) A new sound entry (M_Stamp -> .../ui/Text_Stamp.mp3) is inserted as a
#    new row right before the M_Stage01 row, pushing the bgm rows
#    (M_Stage01..M_Lobby) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Normalize existing paths in column B (rows 2-26) to backslash form ---
$dataRange = $ws.Range("B2:B26")
$dataRange.Replace("../../Resource/sounds", "\..\Resource\sounds")
$dataRange.Replace("/", "\")

# --- 2) Insert the new M_Stamp row above the M_Stage01 row (row 22) ---
$ws.Rows.Item(22).Insert()
$ws.Range("B22").Value = "\..\Resource\sounds\ui\Text_Stamp.mp3"
$ws.Range("A22").Value = "M_Stamp"
$ws.Range("C22").Value = "FMOD_DEFAULT | FMOD_LOOP_OFF"

# --- Update the active selection to reflect the post-edit workbook state ---
$ws.Range("A22").Select()
